# Auto-generated edit script: update cryptos list (Price & Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.923.37"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "1.632.87"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "214.72"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").Value = "0.518"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "28.75"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "1.866.26"
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("D13").Value = "1.639.01"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").Value = "9.29"
$ws.Range("E15").Value = "  +15.02%  "
$ws.Range("D16").Value = "29.935.44"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").Value = "64.15"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "240.77"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "9.82"
$ws.Range("E22").Value = "  +3.40%  "
$ws.Range("E23").Value = "  +2.39%  "
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("D25").Value = "157.92"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").Value = "15.49"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("E31").Value = "  +3.43%  "
$ws.Range("D32").Value = "3.39"
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("D33").Value = "3.18"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("D34").Value = "1.429.86"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("E35").Value = "  +5.56%  "
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").Value = "75.79"
$ws.Range("E40").Value = "  +12.30%  "
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").Value = "1.99"
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("E43").Value = "  +1.25%  "
$ws.Range("D44").Value = "0.0493"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").Value = "51.19"
$ws.Range("E47").Value = "  -7.39%  "
$ws.Range("D48").Value = "5.35"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "1.773.63"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("E50").Value = "  +11.54%  "
$ws.Range("D51").Value = "90.22"
$ws.Range("E51").Value = "  +4.13%  "

